# The original sheet layout was:
#   Row1 headers in B1:F1  -> EL_Astral_exact25, FNRATE_ASTRAL, TAXON, MODEL_CONDITION, GENE
#   Data rows 2-6: column A held a (bold-styled) copy of the GENE number,
#   columns B:F held TAXON, 0, "11-texon", "simulated_25genes_strongILS", GENE number again.
#
# The target layout drops the redundant/duplicated leading column A entirely
# (everything shifts one column to the left: old B->A, C->B, D->C, E->D, F->E),
# giving a clean A1:E6 table, and also fixes a typo in the header text
# "MODEL_CONDITION" -> "MODELCONDITION".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the old column A (the duplicated/styled GENE column); this shifts
# every other column left by one and shrinks the sheet dimension to A1:E6.
$ws.Columns.Item(1).Delete()

# Fix the header typo: MODEL_CONDITION -> MODELCONDITION (now located in D1
# after the column shift above).
$ws.Range("D1").Value = "MODELCONDITION"
